# Preliminary Spoken English Syllabus — populate Sheet1 column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are entered bottom-up (A5 first, A1 last) so the shared-string
# table ends up built in the same order as the authored workbook.
$ws.Range("A5").Value = "WEEK 4- FRAMING OF SENTENCES"
$ws.Range("A4").Value = "WEEK 3- PREPOSITION"
$ws.Range("A3").Value = "WEEK 2- READING AND COMPREHENSION"
$ws.Range("A2").Value = "WEEK 1- EMPHASIS ON PRONOUNCIATION"
$ws.Range("A1").Value = "PRELIMINARY SPOKEN ENGLISH SYLLABUS FOR VICTOR AND SATHI"

# Column A was widened to fit the longest line of text.
$ws.Columns.Item(1).EntireColumn.AutoFit()

# Leave the same cell selected as in the saved workbook.
$ws.Range("E11").Select() | Out-Null
